$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new row for Recife, Brazil (REC) needs to be inserted right before the
# existing "COK" (Kochi, India) row, which is currently row 286. That pushes
# COK down to row 287 and SFO (currently row 287) down to row 288.
#
# Work from the bottom up over just the A:G columns so we don't disturb
# formatting/dimensions outside the data table.

# Shift SFO (287 -> 288), carrying over the colo-code cell's bold/border
# style (column A) along with the values.
$ws.Range("A287").Copy()
$ws.Range("A288").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A288:G288").Value2 = $ws.Range("A287:G287").Value2

# Shift COK (286 -> 287)
$ws.Range("A286").Copy()
$ws.Range("A287").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A287:G287").Value2 = $ws.Range("A286:G286").Value2

# The new Recife row reuses row 286, which already carries the colo-code
# style (it used to hold COK), so no extra formatting step is required here.

# Write the new Recife, Brazil row
$ws.Cells.Item(286, 1).Value = "REC"
$ws.Cells.Item(286, 2).Value = "Recife, Brazil"
$ws.Cells.Item(286, 3).Value = -8.126489639300001
$ws.Cells.Item(286, 4).Value = -34.9235992432
$ws.Cells.Item(286, 5).Value = "BR"
$ws.Cells.Item(286, 6).Value = "South America"
$ws.Cells.Item(286, 7).Value = "Recife"
